$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 4423.75
$ws.Range("J17").Value = 4423.75
$ws.Range("L17").Value = 13271.25
$ws.Range("N17").Value = -13607.25
# Row 21
$ws.Range("H21").Value = 13332.667
$ws.Range("I21").Value = 13332.667
$ws.Range("K21").Value = 13332.667
$ws.Range("M21").Value = -12864.667
# Row 23
$ws.Range("H23").Value = 13332.667
$ws.Range("I23").Value = 13332.667
$ws.Range("K23").Value = 13332.667
$ws.Range("M23").Value = -13098.667
# Row 75
$ws.Range("H75").Value = 69219.336
$ws.Range("J75").Value = 69219.336
$ws.Range("L75").Value = 69219.336
$ws.Range("N75").Value = -71091.336
# Row 78
$ws.Range("H78").Value = 69219.336
$ws.Range("J78").Value = 69219.336
$ws.Range("L78").Value = 207658.008
$ws.Range("N78").Value = -217018.008
# Row 98
$ws.Range("H98").Value = 4000.25
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").Value = $null
# Row 122
$ws.Range("H122").Value = 4000.25
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null
# Row 131
$ws.Range("H131").Value = 1978.909
$ws.Range("I131").Value = 926.8
$ws.Range("J131").Value = 12500
$ws.Range("K131").Value = 2780.4
$ws.Range("L131").Value = 37500
$ws.Range("M131").Value = 2259.6
$ws.Range("N131").Value = -47580

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4509.05
$ws.Range("I61").Value = 1958.7142
$ws.Range("K61").Value = 1958.7142
$ws.Range("M61").Value = -1746.7142
# Row 74
$ws.Range("H74").Value = 2444.15
$ws.Range("I74").Value = 1912.7142
$ws.Range("J74").Value = 3684.1667
$ws.Range("K74").Value = 1912.7142
$ws.Range("L74").Value = 3684.1667
$ws.Range("M74").Value = -1038.7142
$ws.Range("N74").Value = -5432.1667
# Row 77
$ws.Range("H77").Value = 2444.15
$ws.Range("I77").Value = 1912.7142
$ws.Range("J77").Value = 3684.1667
$ws.Range("K77").Value = 9563.571
$ws.Range("L77").Value = 18420.8335
$ws.Range("M77").Value = -5195.571
$ws.Range("N77").Value = -27156.8335
# Row 92
$ws.Range("H92").Value = 61962.332
$ws.Range("J92").Value = 61962.332
$ws.Range("L92").Value = 61962.332
$ws.Range("N92").Value = -66954.33199999999
# Row 136
$ws.Range("H136").Value = 4509.05
$ws.Range("I136").Value = 1958.7142
$ws.Range("K136").Value = 5876.142599999999
$ws.Range("M136").Value = -3326.142599999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 6
$ws.Range("H6").Value = 29247
$ws.Range("J6").Value = 29662.666
$ws.Range("L6").Value = 29662.666
$ws.Range("N6").Value = -29888.666
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = $null
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = $null
# Row 114
$ws.Range("H114").Value = 32221.5
$ws.Range("J114").Value = 32221.5
$ws.Range("L114").Value = 32221.5
$ws.Range("N114").Value = -40899.5
# Row 134
$ws.Range("H134").Value = 1609.75
$ws.Range("I134").Value = 1609.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4829.25
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2294.25
$ws.Range("N134").Value = $null
# Row 92
$ws.Range("H92").Value = 42297.25
$ws.Range("J92").Value = 42297.25
$ws.Range("L92").Value = 42297.25
$ws.Range("N92").Value = -47289.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 1797.5
$ws.Range("I99").Value = 1797.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1797.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -299.5
$ws.Range("N99").Value = $null
# Row 126
$ws.Range("H126").Value = 1797.5
$ws.Range("I126").Value = 1797.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5392.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2922.5
$ws.Range("N126").Value = $null
# Row 132
$ws.Range("H132").Value = 1395.2424
$ws.Range("I132").Value = 1254.7667
$ws.Range("K132").Value = 3764.300099999999
$ws.Range("M132").Value = -1234.300099999999
# Row 134
$ws.Range("H134").Value = 1680.7
$ws.Range("I134").Value = 1610.1111
$ws.Range("K134").Value = 4830.3333
$ws.Range("M134").Value = -2295.3333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Range("H81").Value = 2000
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
# Row 84
$ws.Range("H84").Value = 2000
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
# Row 122
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null
# Row 131
$ws.Range("H131").Value = 2897.1428
$ws.Range("J131").Value = 2897.1428
$ws.Range("L131").Value = 8691.428400000001
$ws.Range("N131").Value = -18771.4284
# Row 139
$ws.Range("H139").Value = 1854.5714
$ws.Range("I139").Value = 1854.5714
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 5563.7142
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -423.7142000000003
$ws.Range("N139").Value = $null
# Row 140
$ws.Range("H140").Value = 5238.3335
$ws.Range("I140").Value = 1333.8667
$ws.Range("K140").Value = 4001.6001
$ws.Range("M140").Value = 1178.3999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 1664875
$ws.Range("J11").Value = 23000
$ws.Range("L11").Value = 23000
$ws.Range("N11").Value = -23278
# Row 70
$ws.Range("H70").Value = 8574.25
$ws.Range("I70").Value = 8725.714
$ws.Range("K70").Value = 8725.714
$ws.Range("M70").Value = -8455.714
# Row 73
$ws.Range("H73").Value = 8574.25
$ws.Range("I73").Value = 8725.714
$ws.Range("K73").Value = 8725.714
$ws.Range("M73").Value = -7789.714
# Row 97
$ws.Range("H97").Value = 225
$ws.Range("I97").Value = 300
$ws.Range("J97").Value = 150
$ws.Range("K97").Value = 300
$ws.Range("L97").Value = 150
$ws.Range("M97").Value = 196
$ws.Range("N97").Value = -1142
# Row 113
$ws.Range("H113").Value = 1168.125
$ws.Range("I113").Value = 1168.125
$ws.Range("K113").Value = 1168.125
$ws.Range("M113").Value = 1001.875
# Row 122
$ws.Range("H122").Value = 3993.158
$ws.Range("I122").Value = 2470.5454
$ws.Range("J122").Value = 6086.75
$ws.Range("K122").Value = 7411.6362
$ws.Range("L122").Value = 18260.25
$ws.Range("M122").Value = -4961.6362
$ws.Range("N122").Value = -23160.25
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1173.5385
$ws.Range("J16").Value = 1783.3334
$ws.Range("L16").Value = 1783.3334
$ws.Range("N16").Value = -2123.3334
# Row 40
$ws.Range("H40").Value = 4037
$ws.Range("I40").Value = 3265.6667
$ws.Range("K40").Value = 3265.6667
$ws.Range("M40").Value = -3129.6667
# Row 68
$ws.Range("H68").Value = 2688
$ws.Range("I68").Value = 2166.3333
$ws.Range("J68").Value = 3001
$ws.Range("K68").Value = 2166.3333
$ws.Range("L68").Value = 3001
$ws.Range("M68").Value = -1417.3333
$ws.Range("N68").Value = -4499
# Row 71
$ws.Range("H71").Value = 2688
$ws.Range("I71").Value = 2166.3333
$ws.Range("J71").Value = 3001
$ws.Range("K71").Value = 10831.6665
$ws.Range("L71").Value = 15005
$ws.Range("M71").Value = -7087.666499999999
$ws.Range("N71").Value = -22493

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 14992.667
$ws.Range("J62").Value = 14994.2
$ws.Range("L62").Value = 14994.2
$ws.Range("N62").Value = -16242.2
# Row 65
$ws.Range("H65").Value = 14992.667
$ws.Range("J65").Value = 14994.2
$ws.Range("L65").Value = 74971
$ws.Range("N65").Value = -81211
# Row 81
$ws.Range("H81").Value = 1667706.5
$ws.Range("I81").Value = 1359.5
$ws.Range("K81").Value = 2719
$ws.Range("M81").Value = -1658
# Row 84
$ws.Range("H84").Value = 1667706.5
$ws.Range("I84").Value = 1359.5
$ws.Range("K84").Value = 13595
$ws.Range("M84").Value = -8291
# Row 100
$ws.Range("H100").Value = 11113880
$ws.Range("I100").Value = 12501865
$ws.Range("K100").Value = 25003730
$ws.Range("M100").Value = -25003189
# Row 107
$ws.Range("H107").Value = 850
$ws.Range("I107").Value = 466.66666
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1399.99998
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = 520.0000199999999
$ws.Range("N107").Value = -9840
# Row 122
$ws.Range("H122").Value = 2000
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900
# Row 132
$ws.Range("H132").Value = 1080.2
$ws.Range("J132").Value = 1401.5
$ws.Range("L132").Value = 4204.5
$ws.Range("N132").Value = -9264.5
# Row 136
$ws.Range("H136").Value = 3286.125
$ws.Range("I136").Value = 3214.8333
$ws.Range("K136").Value = 9644.499899999999
$ws.Range("M136").Value = -7094.499899999999
